$d = $word.ActiveDocument
$para = $d.Paragraphs.First
$para.Range.LanguageID = "en-CA"
$para.Range.Text = "Validate document body content"
$para.Range.LanguageID = "en-CA"
